# --- Build REST API documentation workbook restructuring ---
# 1) Rename the original sheet to "Documentation" and create two more
#    sheets ("For html", "tables") after it, in that order.
$wb = $excel.ActiveWorkbook

$docSheet = $wb.Worksheets.Item(1)
$docSheet.Name = "Documentation"

$htmlSheet = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets.Item("Documentation"))
$htmlSheet.Name = "For html"

$tablesSheet = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets.Item("For html"))
$tablesSheet.Name = "tables"

# Re-fetch fresh (stable) references now that all sheets exist.
$doc = $wb.Worksheets.Item("Documentation")
$html = $wb.Worksheets.Item("For html")
$tables = $wb.Worksheets.Item("tables")

# --- 2) "For html" sheet: copy the request/response summary table (A1:E10)
#    from Documentation, then touch up the two JSON-sample cells which are
#    shorter / unwrapped on this sheet. ---
$doc.Range("A1:E10").Copy($html.Range("A1"))
$html.Range("E3").Value = "JSON object holding data about the movie that was added, including an ID:"
$html.Range("E3").WrapText = $false
$html.Range("E5").Value = "A JSON object holding data about the director that was added, including an ID"
$html.Range("E5").WrapText = $false
$html.Columns.Item(5).ColumnWidth = 67.66666666666667

# --- 3) "tables" sheet: copy the three schema tables that used to live
#    below the main table on Documentation (old B18:G32), shifted up-left
#    by one column and 15 rows (new A3:F17). ---
$doc.Range("B18:G32").Copy($tables.Range("A3"))

# Fix up the handful of cells whose value/style genuinely changed:
$tables.Range("D9").Value = "deathdate"
$tables.Range("E4").Style = "Note"
$tables.Range("E4").Font.Bold = $true
$tables.Range("F4").Style = "Note"
$tables.Range("F4").Font.Bold = $true

$tables.Columns.Item(1).ColumnWidth = 11.666666666666666
$tables.Columns.Item(2).ColumnWidth = 15
$tables.Columns.Item(3).ColumnWidth = 10.166666666666666
$tables.Columns.Item(4).ColumnWidth = 11.166666666666666
$tables.Columns.Item(5).ColumnWidth = 8

$tables.Range("G9").Select()

# --- 4) Documentation sheet: remove the now-relocated schema tables
#    (old rows 18-32) and leave a single "featured" header label behind
#    in G19, matching the new, much wider column G. ---
$doc.Range("B18:G32").Clear()
$doc.Range("G19").Value = "featured"
$doc.Range("G19").Style = "Note"
$doc.Range("G19").Font.Bold = $true

$doc.Columns.Item(7).ColumnWidth = 67.66666666666667

$doc.Range("B3").Select()

# --- 5) Make sure the Documentation tab is the one shown/selected when
#    the workbook re-opens. ---
$doc.Activate()
